$d = $word.ActiveDocument

# --- 0. Remove the old "_GoBack" bookmark first --------------------------
# It currently wraps an empty paragraph right before "Project GitHub…" (its
# matching end tag sits right after the hyperlink). It is going to be
# recreated further down inside the new "Implemented the CGetDialogue…"
# paragraph, so get rid of the old one first to avoid two bookmarks sharing
# the same name at once. Delete() removes both the start and end markers.
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks.Item("_GoBack").Delete()
}

# --- Locate the specific "Finalise UI. (2 hours)" paragraph to edit -------
# There are two identical "Finalise UI. (2 hours)" bullet paragraphs in the
# document; the one we need to change is the second one - the one
# immediately followed by the "Jack Fisher" bullet.
$count = $d.Paragraphs.Count
$targetIndex = $null
for ($i = 1; $i -le $count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Text -eq "Finalise UI. (2 hours)`r") {
        $next = $d.Paragraphs.Item($i + 1)
        if ($next.Range.Text -like "Jack Fisher*") {
            $targetIndex = $i
        }
    }
}

if ($targetIndex -eq $null) {
    Write-Output "ERROR: target paragraph not found"
} else {
    $target = $d.Paragraphs.Item($targetIndex)

    # --- 1. Replace the paragraph's text ------------------------------------
    $textRange = $d.Range($target.Range.Start, $target.Range.End - 1)
    $textRange.Text = "Moved in-app dialogue into an XML file."

    # --- 2. Insert the three new bullet paragraphs after it -----------------
    $target = $d.Paragraphs.Item($targetIndex)
    $insertionPoint = $d.Range($target.Range.End - 1, $target.Range.End - 1)

    $wns = "xmlns:w='http://schemas.openxmlformats.org/wordprocessingml/2006/main'"

    $newParagraphsXml =
        "<w:p $wns>" +
            "<w:pPr>" +
                "<w:pStyle w:val='ListParagraph'/>" +
                "<w:numPr><w:ilvl w:val='1'/><w:numId w:val='2'/></w:numPr>" +
                "<w:spacing w:after='0' w:line='240' w:lineRule='auto'/>" +
            "</w:pPr>" +
            "<w:r><w:t>Created a class to parse the XML file’s strings.</w:t></w:r>" +
        "</w:p>" +
        "<w:p $wns>" +
            "<w:pPr>" +
                "<w:pStyle w:val='ListParagraph'/>" +
                "<w:numPr><w:ilvl w:val='1'/><w:numId w:val='2'/></w:numPr>" +
                "<w:spacing w:after='0' w:line='240' w:lineRule='auto'/>" +
            "</w:pPr>" +
            "<w:r><w:t xml:space='preserve'>Implemented the </w:t></w:r>" +
            "<w:proofErr w:type='spellStart'/>" +
            "<w:r><w:t>CGetDialogue</w:t></w:r>" +
            "<w:proofErr w:type='spellEnd'/>" +
            "<w:r><w:t xml:space='preserve'> class in the main application.</w:t></w:r>" +
            "<w:bookmarkStart w:id='0' w:name='_GoBack'/>" +
            "<w:bookmarkEnd w:id='0'/>" +
        "</w:p>" +
        "<w:p $wns>" +
            "<w:pPr>" +
                "<w:pStyle w:val='ListParagraph'/>" +
                "<w:numPr><w:ilvl w:val='1'/><w:numId w:val='2'/></w:numPr>" +
                "<w:spacing w:after='0' w:line='240' w:lineRule='auto'/>" +
            "</w:pPr>" +
            "<w:r><w:t>Finalise UI. (</w:t></w:r>" +
            "<w:r><w:t>6</w:t></w:r>" +
            "<w:r><w:t xml:space='preserve'> hours)</w:t></w:r>" +
        "</w:p>"

    $insertionPoint.InsertXML($newParagraphsXml) | Out-Null

    Write-Output "Edit applied successfully."
}
